$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to hold purely-numeric-looking text values must be forced
# to Text format first, otherwise Excel will coerce the string into a number.
$numericTextCells = @(
    "H6","H22","H48","H51","H52","H55","H63","H65","H69","H80","H81",
    "G110","H120","H122","H123","H157","H192","G204","G205"
)
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 6
$ws.Range("D6").Value = "2026-02-06T11:20:24.402-0500"
$ws.Range("H6").Value = "24.44"

# Row 22
$ws.Range("D22").Value = "2026-02-06T11:23:45.539-0500"
$ws.Range("H22").Value = "50.65"

# Row 48
$ws.Range("D48").Value = "2026-02-06T10:56:38.686-0500"
$ws.Range("H48").Value = "71.54"

# Row 51
$ws.Range("D51").Value = "2026-02-06T10:55:40.820-0500"
$ws.Range("H51").Value = "72.32"

# Row 52
$ws.Range("D52").Value = "2026-02-06T10:54:16.874-0500"
$ws.Range("H52").Value = "73.19"

# Row 55
$ws.Range("D55").Value = "2026-02-06T10:53:19.380-0500"
$ws.Range("H55").Value = "73.17"

# Row 63
$ws.Range("D63").Value = "2026-02-06T10:51:44.330-0500"
$ws.Range("H63").Value = "73.14"

# Row 65
$ws.Range("D65").Value = "2026-02-06T10:50:41.942-0500"
$ws.Range("H65").Value = "73.13"

# Row 69
$ws.Range("D69").Value = "2026-02-06T10:47:21.620-0500"
$ws.Range("H69").Value = "92.09"

# Row 80
$ws.Range("D80").Value = "2026-02-06T10:44:00.539-0500"
$ws.Range("H80").Value = "97.78"

# Row 81
$ws.Range("D81").Value = "2026-02-06T10:41:40.937-0500"
$ws.Range("H81").Value = "97.74"

# Row 110
$ws.Range("C110").Value = "2026-02-06T15:22:12.479-0500"
$ws.Range("G110").Value = "87.40"

# Row 120
$ws.Range("D120").Value = "2026-02-06T11:25:08.666-0500"
$ws.Range("H120").Value = "172.93"

# Row 122
$ws.Range("D122").Value = "2026-02-06T11:28:02.684-0500"
$ws.Range("H122").Value = "172.98"

# Row 123
$ws.Range("D123").Value = "2026-02-06T11:29:24.363-0500"
$ws.Range("H123").Value = "173.46"

# Row 157
$ws.Range("D157").Value = "2026-02-06T11:31:23.426-0500"
$ws.Range("H157").Value = "242.79"

# Row 192
$ws.Range("D192").Value = "2026-02-06T11:35:00.705-0500"
$ws.Range("H192").Value = "294.73"

# Row 204
$ws.Range("C204").Value = "2026-02-06T15:22:14.589-0500"
$ws.Range("G204").Value = "197.31"

# Row 205
$ws.Range("C205").Value = "2026-02-06T15:22:16.460-0500"
$ws.Range("G205").Value = "197.31"

# Remove the now-obsolete trailing rows 402-410 (data aged out of the report window).
$ws.Range("A402:H410").EntireRow.Delete() | Out-Null
